# Generate Report for Handoff
# Re-run of the localization status report: new source-file UUIDs, a new
# handoff/xliff-generation pass (new hashes + timestamps), and a refreshed
# "Ready for handoff" status (replacing the stale "Handed back" status).

$wb = $excel.ActiveWorkbook

$oldUuid1 = "5950f920-e0d3-4c1e-9147-b86f18c5a9cf"
$oldUuid2 = "dc4c83bd-67b8-4cdd-af8b-a5e8d129ef51"
$newUuid1 = "93b4c6ab-5540-4dd3-b33c-fc9a648f8162"
$newUuid2 = "ffffe88c0b38-6372-4f48-8d62-7a81284afec1"

$oldHash1 = "8a2d1996e1c92ae297805f6ddcd0a3b1086381cc"
$oldHash2 = "8eb2d54313dd4cd85116b56048ddb09644c904e0"
$newHash  = "f9615f04b769825cba53aa39339543ff522f2111"

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newUuid1.md"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-27 06:59:19"

$wsOverview.Range("A3").Value = "$newUuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newUuid2.md"
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = "2016-08-27 06:59:19"

# Refresh the hyperlink display text (targets are unchanged).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$oldUuid1.md", "", "", "e2e\$newUuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$oldUuid2.md", "", "", "e2e\$newUuid2.md")

$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newUuid1.md"
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("G2").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-27 06:59:14"
$wsZh.Range("I2").ClearContents()
$wsZh.Range("J2").ClearContents()
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Range("A3").Value = "$newUuid2.md"
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-27 06:59:14"
$wsZh.Range("I3").ClearContents()
$wsZh.Range("J3").ClearContents()
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"

# Reset the (now plain, non-hyperlinked) I2/I3 cell style back to Normal.
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("I3").Style = "Normal"

# Refresh the remaining hyperlinks (A2/A3); targets are unchanged, only the
# display text changes. The I2/I3 hyperlinks are intentionally dropped.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$oldUuid1.md", "", "", "$newUuid1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$oldUuid2.md", "", "", "$newUuid2.md")

$wsZh.Columns.Item(3).ColumnWidth = 16.33
$wsZh.Columns.Item(9).ColumnWidth = 17.83
$wsZh.Columns.Item(10).ColumnWidth = 20.83

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newUuid1.md"
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("G2").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-27 06:59:19"
$wsDe.Range("I2").ClearContents()
$wsDe.Range("J2").ClearContents()
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Range("A3").Value = "$newUuid2.md"
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-27 06:59:19"
$wsDe.Range("I3").ClearContents()
$wsDe.Range("J3").ClearContents()
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"

# Reset the (now plain, non-hyperlinked) I2/I3 cell style back to Normal.
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("I3").Style = "Normal"

# Refresh the remaining hyperlinks (A2/A3); targets are unchanged, only the
# display text changes. The I2/I3 hyperlinks are intentionally dropped.
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$oldUuid1.md", "", "", "$newUuid1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d04e6250018528f9c027b38135e1f502e7e617/e2e/$oldUuid2.md", "", "", "$newUuid2.md")

$wsDe.Columns.Item(3).ColumnWidth = 16.33
$wsDe.Columns.Item(9).ColumnWidth = 17.83
$wsDe.Columns.Item(10).ColumnWidth = 20.83
